# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 550, shifting the existing
# rows 550-623 down to 551-624 (dimension grows from A1:R623 to A1:R624).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 550; Excel shifts rows 550..623 down to 551..624.
$ws.Rows.Item(550).Insert()

# Populate the newly inserted row 550 with the new record's data.
$ws.Cells.Item(550, 1).Value  = 3
$ws.Cells.Item(550, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(550, 3).Value  = "Coquimbo"
$ws.Cells.Item(550, 4).Value  = 45154
$ws.Cells.Item(550, 5).Value  = 5
$ws.Cells.Item(550, 6).Value  = 100112040
$ws.Cells.Item(550, 7).Value  = "Cilantro"
$ws.Cells.Item(550, 8).Value  = "Sin especificar"
$ws.Cells.Item(550, 9).Value  = "Primera"
$ws.Cells.Item(550, 10).Value = 120
$ws.Cells.Item(550, 11).Value = 4000
$ws.Cells.Item(550, 12).Value = 4000
$ws.Cells.Item(550, 13).Value = 4000
$ws.Cells.Item(550, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(550, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(550, 16).Value = 1333
$ws.Cells.Item(550, 17).Value = 3
$ws.Cells.Item(550, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D (s="2").
$ws.Cells.Item(550, 4).NumberFormat = $ws.Cells.Item(551, 4).NumberFormat
